$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2: Nordsjaelland vs Vejle, 21:00, MS 5.5 Üst, 6.98
$ws.Range("A2").Value = "Nordsjaelland"
$ws.Range("B2").Value = "Vejle"
$ws.Range("C2").Value = "21:00"
$ws.Range("D2").Value = "MS 5.5 Üst"
$ws.Range("E2").Value = 6.98

# Update row 3: Helmond vs Den Haag, 22:00, MS 5.5 Üst, 5.72
$ws.Range("A3").Value = "Helmond"
$ws.Range("B3").Value = "Den Haag"
$ws.Range("C3").Value = "22:00"
$ws.Range("D3").Value = "MS 5.5 Üst"
$ws.Range("E3").Value = 5.72

# Remove old prediction rows 4-10 (clear contents)
$ws.Range("A4:E10").Clear()
